$wb = $excel.ActiveWorkbook

$wsRecipients = $wb.Worksheets.Item("Recipients")
$wsWishes = $wb.Worksheets.Item("Wishes")

# --- Add the four new recipients (rows 7-10) ---
# Order matters: these new shared strings must be appended in this exact
# row-major (A then B) order so the shared-string table ends up matching
# the target workbook.
$wsRecipients.Range("A7").Value2 = "Оптимус Прайм"
$wsRecipients.Range("B7").Value2 = "трансформер"

$wsRecipients.Range("A8").Value2 = "капитан Шепард"
$wsRecipients.Range("B8").Value2 = "СПЕКТР"

$wsRecipients.Range("A9").Value2 = "Трисс Меригольд"
$wsRecipients.Range("B9").Value2 = "жен"

$wsRecipients.Range("A10").Value2 = "Геральт из Ривии"
$wsRecipients.Range("B10").Value2 = "мужчина"

# --- Widen column B on the Recipients sheet to fit the new values ---
$wsRecipients.Columns.Item(2).ColumnWidth = 15

# --- Give the Recipients sheet explicit page setup (paper size & orientation) ---
$wsRecipients.PageSetup.PaperSize = 9
$wsRecipients.PageSetup.Orientation = 1

# --- Update the wish text on the Wishes sheet (row 5, column E) ---
# This replaces the old, now-unused wish string with a new one; the
# old string becomes orphaned and is dropped from the shared string table
# while the new one is appended last.
$wsWishes.Range("E5").Value2 = "исполнения всех желаний"

# --- Update selections on each sheet ---
$wsRecipients.Range("E8").Select()

# --- Make Wishes the active sheet/tab and set its selection ---
$wsWishes.Activate()
$wsWishes.Range("B2").Select()
